$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.649422333333334
$ws.Range("H2").Value = 22.948267
$ws.Range("I2").Value = 0.004484559810904267
$ws.Range("J2").Value = 0.004484559810904268
$ws.Range("M2").Value = 3.456265333333333
$ws.Range("N2").Value = 10.368796
$ws.Range("O2").Value = 0.009841535807677501
$ws.Range("P2").Value = 0.0098415358076775
$ws.Range("Q2").Value = 26.43843323072578
$ws.Range("R2").Value = 237.945899076532
$ws.Range("S2").Value = 0.00004413495596068579
$ws.Range("T2").Value = 0.0000441349559606858
$ws.Range("G3").Value = 7.649422333333334
$ws.Range("H3").Value = 22.948267
$ws.Range("I3").Value = 0.004484559810904267
$ws.Range("J3").Value = 0.004484559810904268
$ws.Range("O3").Value = 0.8587907398420774
$ws.Range("P3").Value = 0.8587907398420773
$ws.Range("Q3").Value = 2307.066912947454
$ws.Range("R3").Value = 20763.60221652708
$ws.Range("S3").Value = 0.003851298437872522
$ws.Range("T3").Value = 0.003851298437872523
$ws.Range("G4").Value = 7.649422333333334
$ws.Range("H4").Value = 22.948267
$ws.Range("I4").Value = 0.004484559810904267
$ws.Range("J4").Value = 0.004484559810904268
$ws.Range("O4").Value = 0.1313677243502452
$ws.Range("P4").Value = 0.1313677243502452
$ws.Range("Q4").Value = 352.9080091541093
$ws.Range("R4").Value = 3176.172082386984
$ws.Range("S4").Value = 0.0005891264170710595
$ws.Range("T4").Value = 0.0005891264170710595
$ws.Range("I5").Value = 0.8893308176045429
$ws.Range("J5").Value = 0.889330817604543
$ws.Range("M5").Value = 3.456265333333333
$ws.Range("N5").Value = 10.368796
$ws.Range("O5").Value = 0.009841535807677501
$ws.Range("P5").Value = 0.0098415358076775
$ws.Range("Q5").Value = 5242.992497077077
$ws.Range("R5").Value = 47186.9324736937
$ws.Range("S5").Value = 0.008752381086326219
$ws.Range("T5").Value = 0.008752381086326217
$ws.Range("I6").Value = 0.8893308176045429
$ws.Range("J6").Value = 0.889330817604543
$ws.Range("O6").Value = 0.8587907398420774
$ws.Range("P6").Value = 0.8587907398420773
$ws.Range("S6").Value = 0.763749070814965
$ws.Range("T6").Value = 0.763749070814965
$ws.Range("I7").Value = 0.8893308176045429
$ws.Range("J7").Value = 0.889330817604543
$ws.Range("O7").Value = 0.1313677243502452
$ws.Range("P7").Value = 0.1313677243502452
$ws.Range("S7").Value = 0.1168293657032518
$ws.Range("T7").Value = 0.1168293657032518
$ws.Range("I8").Value = 0.1061846225845528
$ws.Range("J8").Value = 0.1061846225845528
$ws.Range("M8").Value = 3.456265333333333
$ws.Range("N8").Value = 10.368796
$ws.Range("O8").Value = 0.009841535807677501
$ws.Range("P8").Value = 0.0098415358076775
$ws.Range("Q8").Value = 626.0045963720661
$ws.Range("R8").Value = 5634.041367348595
$ws.Range("S8").Value = 0.001045019765390597
$ws.Range("T8").Value = 0.001045019765390597
$ws.Range("I9").Value = 0.1061846225845528
$ws.Range("J9").Value = 0.1061846225845528
$ws.Range("O9").Value = 0.8587907398420774
$ws.Range("P9").Value = 0.8587907398420773
$ws.Range("S9").Value = 0.09119037058923986
$ws.Range("T9").Value = 0.09119037058923986
$ws.Range("I10").Value = 0.1061846225845528
$ws.Range("J10").Value = 0.1061846225845528
$ws.Range("O10").Value = 0.1313677243502452
$ws.Range("P10").Value = 0.1313677243502452
$ws.Range("S10").Value = 0.01394923222992235
$ws.Range("T10").Value = 0.01394923222992235
